# EUC Perth Assets — Build Room tracker update
#
# 1) "4.2 Items": the Desktop Mini G9 tally is corrected — LastCount moves
#    63 -> 64 and NewCount moves 64 -> 63 (net -1), reflecting the
#    "Subtract 1" action logged below.
# 2) "4.2 Timestamps": the four most-recent log rows (41-44) get the same
#    explicit cell formatting as every earlier row in the log, and a new
#    row 45 is appended recording a "Subtract 1" for the Desktop Mini G9.

$wb = $excel.ActiveWorkbook

# --- 1) "4.2 Items": swap the Desktop Mini G9 LastCount / NewCount ---
$items = $wb.Worksheets.Item("4.2 Items")
$items.Range("B2").Value = 64
$items.Range("C2").Value = 63

# --- 2) "4.2 Timestamps": normalize formatting on rows 41-44, append row 45 ---
$log = $wb.Worksheets.Item("4.2 Timestamps")

# Every earlier row (1-40) in the log carries an explicit cell style; rows
# 41-44 were never normalized. Bring them into line using the style already
# applied to row 40.
$normalStyle = $log.Range("A40").Style
$log.Range("A41:D44").Style = $normalStyle

# Append the new log entry for row 45 (left unstyled, same as 41-44 were
# before this pass, since it is the freshest/raw entry).
$log.Range("A45").Value = "Desktop Mini G9"
$log.Range("B45").Value = "Subtract 1"
$log.Range("C45").Value = ""
$log.Range("D45").Value = "2023-12-28 14:01:48"
